$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "D"

$ws.Range("B2").Value = "1"
$ws.Range("D2").Value = "1"

$ws.Range("C3").Value = "1"

$ws.Range("D4").Value = "1"
